# The commit renames the display "Name" of three inline logo pictures
# that live in the document's headers/footers:
#   - the Pearson logo inline picture (currently "image2.png") -> "image1.png"
#     (it occurs in both the "first page" footer and the "default" footer)
#   - the BTEC logo inline picture (currently "image1.jpg") -> "image2.jpg"
#     (it occurs in the "first page" header)
# Nothing else about the pictures (their embedded binary, size, description,
# position, ...) changes - only the <wp:docPr>/<pic:cNvPr> "name" attribute.
#
# NOTE: InlineShape.Name is a write-only-ish property here (reading it back
# before ever assigning it does not reflect the stored docPr/name), so shapes
# are identified by their (stable, readable) AlternativeText / description
# instead of by their current Name.

$d = $word.ActiveDocument

function Rename-InlineShapesInRange($range) {
    if ($range -eq $null) { return }
    $shapes = $range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $descr = $shp.AlternativeText

        if ($descr -eq "BTec_Logo-Orange") {
            $shp.Name = "image2.jpg"
        } elseif ($descr -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image1.png"
        }
    }
}

# Inline pictures that live in the main document body (none in this file,
# but handled for completeness/robustness).
Rename-InlineShapesInRange $d.Content

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2, wdHeaderFooterEvenPages = 3
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    for ($h = 1; $h -le $section.Headers.Count; $h++) {
        $header = $section.Headers.Item($h)
        if ($header.Exists) {
            Rename-InlineShapesInRange $header.Range
        }
    }

    for ($f = 1; $f -le $section.Footers.Count; $f++) {
        $footer = $section.Footers.Item($f)
        if ($footer.Exists) {
            Rename-InlineShapesInRange $footer.Range
        }
    }
}
